$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.461107730865479
$ws.Range("B1").Value = 2.940497636795044
$ws.Range("C1").Value = 2.585555791854858
$ws.Range("D1").Value = 2.385602474212646
$ws.Range("E1").Value = 1.742071151733398
